# Add 2022-Q3 data
# 1) Insert a new row into the "总计" (summary) sheet for the 2022-Q3 quarter.
# 2) Insert a brand-new "2022-Q3" worksheet (cloned from the existing "2022-Q1"
#    sheet so that formatting/styles match exactly) right after "总计", holding
#    the underlying fund holdings for that quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" sheet: insert new row 2 with the 2022-Q3 summary figures.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Rows.Item(2).Insert()
$summary.Range("A2:D2").ClearFormats()

# Match the existing index-column styling (same as A3, A4, ...).
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 0.08

# Re-number the row index column so it stays a simple 0..6 sequence.
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
$summary.Range("A7").Value = 5
$summary.Range("A8").Value = 6

# ---------------------------------------------------------------------------
# 2) New "2022-Q3" worksheet, placed right after "总计" and before "2022-Q1".
# ---------------------------------------------------------------------------
$quarterOne = $wb.Worksheets.Item(2)
$quarterOne.Copy($quarterOne, $null)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q3"

# Drop the two extra rows that "2022-Q1" had (rows 4 and 5) — 2022-Q3 only
# has two fund rows.
$newSheet.Range("A4:H5").Clear()

# Fund-code / numeric-as-text columns need to keep their text formatting so
# leading zeros and fixed decimals are preserved (same as the other sheets).
$newSheet.Range("B2:B3").NumberFormat = "@"
$newSheet.Range("D2:G3").NumberFormat = "@"

$newSheet.Range("B2").Value = "005052"
$newSheet.Range("C2").Value = "上投摩根标普港股通低波红利指数C"
$newSheet.Range("D2").Value = "1.37"
$newSheet.Range("E2").Value = "92.94"
$newSheet.Range("F2").Value = "2.78"
$newSheet.Range("G2").Value = "0.0381"
$newSheet.Range("H2").Value = 3

$newSheet.Range("B3").Value = "005051"
$newSheet.Range("C3").Value = "上投摩根标普港股通低波红利指数A"
$newSheet.Range("D3").Value = "1.36"
$newSheet.Range("E3").Value = "92.94"
$newSheet.Range("F3").Value = "2.78"
$newSheet.Range("G3").Value = "0.0378"
$newSheet.Range("H3").Value = 3
